# Applies the row updates for "Paraguay Division Profesional.xlsx" -
# re-ordering of match records (ids/odds) within matching dates, plus
# a brand-new match recorded in what was previously a partial row 185.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 134
$ws.Cells.Item(134, 2).Value = 7493428
$ws.Cells.Item(134, 6).Value = "Guairena FC"
$ws.Cells.Item(134, 7).Value = "Resistencia FC"
$ws.Cells.Item(134, 8).Value = 4
$ws.Cells.Item(134, 9).Value = 1
$ws.Cells.Item(134, 10).Value = "H"
$ws.Cells.Item(134, 11).Value = 1.727
$ws.Cells.Item(134, 12).Value = 3.6
$ws.Cells.Item(134, 13).Value = 4.2
$ws.Cells.Item(134, 14).Value = 1.45
$ws.Cells.Item(134, 15).Value = 4.2
$ws.Cells.Item(134, 16).Value = 6
$ws.Cells.Item(134, 17).Value = -1
$ws.Cells.Item(134, 18).Value = 1.775
$ws.Cells.Item(134, 19).Value = 2.025
$ws.Cells.Item(134, 20).Value = 2.75
$ws.Cells.Item(134, 21).Value = 1.825
$ws.Cells.Item(134, 22).Value = 1.975
$ws.Cells.Item(134, 23).Value = 0.45
$ws.Cells.Item(134, 24).Value = -1
$ws.Cells.Item(134, 25).Value = -1
$ws.Cells.Item(134, 26).Value = 0.7749999999999999
$ws.Cells.Item(134, 27).Value = -1
$ws.Cells.Item(134, 28).Value = 0.825
$ws.Cells.Item(134, 29).Value = -1

# Row 135
$ws.Cells.Item(135, 2).Value = 7493427
$ws.Cells.Item(135, 6).Value = "Tacuary"
$ws.Cells.Item(135, 7).Value = "Sportivo Luqueno"
$ws.Cells.Item(135, 8).Value = 1
$ws.Cells.Item(135, 9).Value = 1
$ws.Cells.Item(135, 10).Value = "D"
$ws.Cells.Item(135, 11).Value = 3.4
$ws.Cells.Item(135, 12).Value = 3.3
$ws.Cells.Item(135, 13).Value = 2
$ws.Cells.Item(135, 14).Value = 3.2
$ws.Cells.Item(135, 15).Value = 3.25
$ws.Cells.Item(135, 16).Value = 2.1
$ws.Cells.Item(135, 17).Value = 0.25
$ws.Cells.Item(135, 18).Value = 2.025
$ws.Cells.Item(135, 19).Value = 1.775
$ws.Cells.Item(135, 20).Value = 2.5
$ws.Cells.Item(135, 21).Value = 1.975
$ws.Cells.Item(135, 22).Value = 1.825
$ws.Cells.Item(135, 23).Value = -1
$ws.Cells.Item(135, 24).Value = 2.25
$ws.Cells.Item(135, 25).Value = -1
$ws.Cells.Item(135, 26).Value = 0.5125
$ws.Cells.Item(135, 27).Value = -0.5
$ws.Cells.Item(135, 28).Value = -1
$ws.Cells.Item(135, 29).Value = 0.825

# Row 140
$ws.Cells.Item(140, 2).Value = 7493431
$ws.Cells.Item(140, 6).Value = "Sportivo Trinidense"
$ws.Cells.Item(140, 7).Value = "Guairena FC"
$ws.Cells.Item(140, 8).Value = 7
$ws.Cells.Item(140, 9).Value = 2
$ws.Cells.Item(140, 10).Value = "H"
$ws.Cells.Item(140, 11).Value = 2.05
$ws.Cells.Item(140, 12).Value = 3.3
$ws.Cells.Item(140, 13).Value = 3.3
$ws.Cells.Item(140, 14).Value = 2.6
$ws.Cells.Item(140, 15).Value = 3.1
$ws.Cells.Item(140, 16).Value = 2.6
$ws.Cells.Item(140, 17).Value = 0
$ws.Cells.Item(140, 18).Value = 1.925
$ws.Cells.Item(140, 19).Value = 1.875
$ws.Cells.Item(140, 20).Value = 2.5
$ws.Cells.Item(140, 21).Value = 2
$ws.Cells.Item(140, 22).Value = 1.8
$ws.Cells.Item(140, 23).Value = 1.6
$ws.Cells.Item(140, 24).Value = -1
$ws.Cells.Item(140, 25).Value = -1
$ws.Cells.Item(140, 26).Value = 0.925
$ws.Cells.Item(140, 27).Value = -1
$ws.Cells.Item(140, 28).Value = 1
$ws.Cells.Item(140, 29).Value = -1

# Row 141
$ws.Cells.Item(141, 2).Value = 7493310
$ws.Cells.Item(141, 6).Value = "Libertad Asuncion"
$ws.Cells.Item(141, 7).Value = "Tacuary"
$ws.Cells.Item(141, 8).Value = 1
$ws.Cells.Item(141, 9).Value = 2
$ws.Cells.Item(141, 10).Value = "A"
$ws.Cells.Item(141, 11).Value = 1.363
$ws.Cells.Item(141, 12).Value = 5
$ws.Cells.Item(141, 13).Value = 7
$ws.Cells.Item(141, 14).Value = 1.571
$ws.Cells.Item(141, 15).Value = 4.2
$ws.Cells.Item(141, 16).Value = 4.75
$ws.Cells.Item(141, 17).Value = -0.75
$ws.Cells.Item(141, 18).Value = 1.8
$ws.Cells.Item(141, 19).Value = 2
$ws.Cells.Item(141, 20).Value = 2.75
$ws.Cells.Item(141, 21).Value = 1.8
$ws.Cells.Item(141, 22).Value = 2
$ws.Cells.Item(141, 23).Value = -1
$ws.Cells.Item(141, 24).Value = -1
$ws.Cells.Item(141, 25).Value = 3.75
$ws.Cells.Item(141, 26).Value = -1
$ws.Cells.Item(141, 27).Value = 1
$ws.Cells.Item(141, 28).Value = 0.4
$ws.Cells.Item(141, 29).Value = -0.5

# Row 143
$ws.Cells.Item(143, 2).Value = 7493311
$ws.Cells.Item(143, 6).Value = "General Caballero JLM"
$ws.Cells.Item(143, 7).Value = "Olimpia Asuncion"
$ws.Cells.Item(143, 8).Value = 0
$ws.Cells.Item(143, 9).Value = 1
$ws.Cells.Item(143, 10).Value = "A"
$ws.Cells.Item(143, 11).Value = 3.4
$ws.Cells.Item(143, 12).Value = 3.3
$ws.Cells.Item(143, 13).Value = 2
$ws.Cells.Item(143, 14).Value = 3.2
$ws.Cells.Item(143, 15).Value = 3.25
$ws.Cells.Item(143, 16).Value = 2.1
$ws.Cells.Item(143, 17).Value = 0.25
$ws.Cells.Item(143, 18).Value = 1.95
$ws.Cells.Item(143, 19).Value = 1.85
$ws.Cells.Item(143, 20).Value = 2.25
$ws.Cells.Item(143, 21).Value = 1.775
$ws.Cells.Item(143, 22).Value = 2.025
$ws.Cells.Item(143, 23).Value = -1
$ws.Cells.Item(143, 24).Value = -1
$ws.Cells.Item(143, 25).Value = 1.1
$ws.Cells.Item(143, 26).Value = -1
$ws.Cells.Item(143, 27).Value = 0.8500000000000001
$ws.Cells.Item(143, 28).Value = -1
$ws.Cells.Item(143, 29).Value = 1.025

# Row 144
$ws.Cells.Item(144, 2).Value = 7493312
$ws.Cells.Item(144, 6).Value = "Cerro Porteno"
$ws.Cells.Item(144, 7).Value = "Guarani Asuncion"
$ws.Cells.Item(144, 8).Value = 4
$ws.Cells.Item(144, 9).Value = 0
$ws.Cells.Item(144, 10).Value = "H"
$ws.Cells.Item(144, 11).Value = 1.7
$ws.Cells.Item(144, 12).Value = 3.6
$ws.Cells.Item(144, 13).Value = 4.333
$ws.Cells.Item(144, 14).Value = 1.727
$ws.Cells.Item(144, 15).Value = 3.75
$ws.Cells.Item(144, 16).Value = 4.2
$ws.Cells.Item(144, 17).Value = -0.5
$ws.Cells.Item(144, 18).Value = 1.8
$ws.Cells.Item(144, 19).Value = 2
$ws.Cells.Item(144, 20).Value = 2.75
$ws.Cells.Item(144, 21).Value = 1.875
$ws.Cells.Item(144, 22).Value = 1.925
$ws.Cells.Item(144, 23).Value = 0.7270000000000001
$ws.Cells.Item(144, 24).Value = -1
$ws.Cells.Item(144, 25).Value = -1
$ws.Cells.Item(144, 26).Value = 0.8
$ws.Cells.Item(144, 27).Value = -1
$ws.Cells.Item(144, 28).Value = 0.875
$ws.Cells.Item(144, 29).Value = -1

# Row 145
$ws.Cells.Item(145, 2).Value = 7493433
$ws.Cells.Item(145, 6).Value = "Sportivo Luqueno"
$ws.Cells.Item(145, 7).Value = "Nacional Asuncion"
$ws.Cells.Item(145, 8).Value = 1
$ws.Cells.Item(145, 9).Value = 1
$ws.Cells.Item(145, 10).Value = "D"
$ws.Cells.Item(145, 11).Value = 2.75
$ws.Cells.Item(145, 12).Value = 3.2
$ws.Cells.Item(145, 13).Value = 2.4
$ws.Cells.Item(145, 14).Value = 2.75
$ws.Cells.Item(145, 15).Value = 3.1
$ws.Cells.Item(145, 16).Value = 2.45
$ws.Cells.Item(145, 17).Value = 0.25
$ws.Cells.Item(145, 18).Value = 1.75
$ws.Cells.Item(145, 19).Value = 2.05
$ws.Cells.Item(145, 20).Value = 2.25
$ws.Cells.Item(145, 21).Value = 2
$ws.Cells.Item(145, 22).Value = 1.8
$ws.Cells.Item(145, 23).Value = -1
$ws.Cells.Item(145, 24).Value = 2.1
$ws.Cells.Item(145, 25).Value = -1
$ws.Cells.Item(145, 26).Value = 0.375
$ws.Cells.Item(145, 27).Value = -0.5
$ws.Cells.Item(145, 28).Value = -0.5
$ws.Cells.Item(145, 29).Value = 0.4

# Row 185
$ws.Cells.Item(185, 5).Value = 45347.86458333334
$ws.Cells.Item(185, 2).Value = 7609137
$ws.Cells.Item(185, 6).Value = "Guarani Asuncion"
$ws.Cells.Item(185, 7).Value = "Olimpia Asuncion"
$ws.Cells.Item(185, 8).Value = 1
$ws.Cells.Item(185, 9).Value = 3
$ws.Cells.Item(185, 10).Value = "A"
$ws.Cells.Item(185, 11).Value = 2.6
$ws.Cells.Item(185, 12).Value = 3.2
$ws.Cells.Item(185, 13).Value = 2.5
$ws.Cells.Item(185, 14).Value = 2.8
$ws.Cells.Item(185, 15).Value = 3.2
$ws.Cells.Item(185, 16).Value = 2.3
$ws.Cells.Item(185, 17).Value = 0.25
$ws.Cells.Item(185, 18).Value = 1.775
$ws.Cells.Item(185, 19).Value = 2.025
$ws.Cells.Item(185, 20).Value = 2.25
$ws.Cells.Item(185, 21).Value = 1.9
$ws.Cells.Item(185, 22).Value = 1.9
$ws.Cells.Item(185, 23).Value = -1
$ws.Cells.Item(185, 24).Value = -1
$ws.Cells.Item(185, 25).Value = 1.3
$ws.Cells.Item(185, 26).Value = -1
$ws.Cells.Item(185, 27).Value = 1.025
$ws.Cells.Item(185, 28).Value = 0.8999999999999999
$ws.Cells.Item(185, 29).Value = -1
